# Auto-generated edit script: updates cryptos list values (Price / Volume(1h) / swapped rows)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '40.029.63'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -3.30%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.335.20'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -4.36%  '

# Row 4
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.20%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '308.92'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -2.40%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '85.10'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -5.43%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.529'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -2.65%  '

# Row 8
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -0.02%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.484'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -2.83%  '

# Row 10
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -1.87%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '30.18'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -6.33%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.110'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +1.05%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '2.694.77'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -4.36%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.44'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -4.13%  '

# Row 15
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -4.17%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.340.46'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -3.73%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.758'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -1.84%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '39.984.98'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -3.13%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.0₃0902'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -2.15%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.11'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -2.25%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '67.93'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -5.68%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '10.66'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -4.29%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '236.05'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +0.24%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.55'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -5.45%  '

# Row 25
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +0.17%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.81'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -4.26%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '23.31'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -3.16%  '

# Row 28
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -4.52%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.26'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -3.20%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '34.18'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -1.81%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '153.76'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -1.94%  '

# Row 32
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -0.20%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.11'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -3.08%  '

# Row 34
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -4.19%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0717'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -3.82%  '

# Row 36
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -0.71%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0998'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -0.47%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.74'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -5.75%  '

# Row 39
$ws.Range("B39").Value = 'ARBITRUM'
$ws.Range("C39").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.72'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -2.50%  '

# Row 40
$ws.Range("B40").Value = 'Celestia'
$ws.Range("C40").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '15.55'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -6.28%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.82'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -1.94%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.946.97'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -1.89%  '

# Row 43
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -4.74%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '17.61'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -3.19%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0263'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -4.73%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '9.25'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -2.32%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.71'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -5.53%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.553.99'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -4.73%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '92.69'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -2.88%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '70.65'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -3.84%  '

# Row 51
$ws.Range("B51").Value = 'MultiversX'
$ws.Range("C51").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '50.01'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -4.18%  '
